# aggiornamento 15, 16, 17 marzo
# Appends three new daily rows (227-229) to the single data sheet,
# continuing the existing table of date / nuovi pos. / somma mobile 7gg. /
# somma mobile 7gg. per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 227 - 2021-04-15 (44301)
$ws.Range("A226").Copy()
$ws.Range("A227").PasteSpecial(-4122)   # xlPasteFormats: reuse A226's style (s="2") without touching styles.xml
$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 0
$ws.Range("C227").Value = 1
$ws.Range("D227").Value = 18.93939393939394

# Row 228 - 2021-04-16 (44302)
$ws.Range("A226").Copy()
$ws.Range("A228").PasteSpecial(-4122)
$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 2
$ws.Range("C228").Value = 3
$ws.Range("D228").Value = 56.81818181818181

# Row 229 - 2021-04-17 (44303)
$ws.Range("A226").Copy()
$ws.Range("A229").PasteSpecial(-4122)
$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 0
$ws.Range("C229").Value = 2
$ws.Range("D229").Value = 37.87878787878788

$excel.CutCopyMode = $false
